$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Avocado", "Dark Green"),
    @("Banana", "Yellow"),
    @("Watermelon", "Bottle Green"),
    @("Kiwi", "Brown"),
    @("Black Currant", "Black")
)

$startRow = 8
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    if ($row -eq 12) {
        # Ensure "Black" is registered in the shared strings table before
        # "Black Currant" so the unique string order matches the source.
        $ws.Cells.Item($row, 2).Value = $data[$i][1]
        $ws.Cells.Item($row, 1).Value = $data[$i][0]
    } else {
        $ws.Cells.Item($row, 1).Value = $data[$i][0]
        $ws.Cells.Item($row, 2).Value = $data[$i][1]
    }
}

$ws.Range("A13").Select()
